$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoiFormulaHelperTest")

$ws.Range("A3").Formula = "=(6*12)+2"
$ws.Range("A4").Formula = "=6*(12+2)"

$ws.Range("B7").Select()
